$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column keeps its text formatting (values look numeric
# like "25.388.49" or "1.0000" and must not be auto-converted to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.388.49'
$ws.Range('E2').Value = '  -1.99%  '

$ws.Range('D3').Value = '1.663.82'
$ws.Range('E3').Value = '  -3.93%  '

$ws.Range('D4').Value = '0.9984'
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').Value = '236.04'
$ws.Range('E5').Value = '  -3.80%  '

$ws.Range('D6').Value = '0.9992'
$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('D7').Value = '0.4801'
$ws.Range('E7').Value = '  -4.61%  '

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').Value = '39.76'
$ws.Range('E8').Value = '  -2.08%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '0.2616'
$ws.Range('E9').Value = '  -3.77%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.06145'
$ws.Range('E10').Value = '  -0.42%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.07078'
$ws.Range('E11').Value = '  -2.20%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.658.55'
$ws.Range('E12').Value = '  -4.50%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '14.72'
$ws.Range('E13').Value = '  -3.11%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.5922'
$ws.Range('E14').Value = '  -9.37%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '4.380'
$ws.Range('E15').Value = '  -7.93%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '74.29'
$ws.Range('E16').Value = '  -3.53%  '

$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').Value = '0.9990'
$ws.Range('E17').Value = '  -0.14%  '

$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = '0.9993'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '25.379.22'
$ws.Range('E19').Value = '  -2.04%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.000006702'
$ws.Range('E20').Value = '  -1.86%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '11.41'
$ws.Range('E21').Value = '  -4.17%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '1.869.27'
$ws.Range('E22').Value = '  -4.58%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '4.425'
$ws.Range('E23').Value = '  -3.55%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '8.637'
$ws.Range('E24').Value = '  -1.93%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '5.329'
$ws.Range('E25').Value = '  -2.69%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '133.84'
$ws.Range('E26').Value = '  -0.47%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '15.10'
$ws.Range('E27').Value = '  -1.25%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.396'
$ws.Range('E28').Value = '  -1.74%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '104.69'
$ws.Range('E29').Value = '  -0.67%  '

$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = '1.689'
$ws.Range('E30').Value = '  -5.54%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '3.979'
$ws.Range('E31').Value = '  +0.36%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '3.624'
$ws.Range('E32').Value = '  -2.09%  '

$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').Value = '0.07656'
$ws.Range('E33').Value = '  -5.87%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04386'
$ws.Range('E34').Value = '  -7.15%  '

$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = '0.9983'
$ws.Range('E35').Value = '  -0.09%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.599'
$ws.Range('E36').Value = '  -1.98%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.6034'
$ws.Range('E37').Value = '  -1.59%  '

$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '0.9421'
$ws.Range('E38').Value = '  -5.31%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.622'
$ws.Range('E39').Value = '  -4.51%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '0.8522'
$ws.Range('E40').Value = '  -3.17%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.0000'
$ws.Range('E41').Value = '  -0.02%  '

$ws.Range('D42').Value = '0.01502'
$ws.Range('E42').Value = '  -6.69%  '

$ws.Range('D43').Value = '1.826'
$ws.Range('E43').Value = '  -6.86%  '

$ws.Range('D44').Value = '98.57'
$ws.Range('E44').Value = '  -3.16%  '

$ws.Range('D45').Value = '0.3756'
$ws.Range('E45').Value = '  -3.87%  '

$ws.Range('D46').Value = '4.704'
$ws.Range('E46').Value = '  -6.14%  '

$ws.Range('D47').Value = '0.1115'
$ws.Range('E47').Value = '  -5.86%  '

$ws.Range('D48').Value = '6.210'
$ws.Range('E48').Value = '  -2.49%  '

$ws.Range('D49').Value = '0.05247'
$ws.Range('E49').Value = '  -0.61%  '

$ws.Range('D50').Value = '29.48'
$ws.Range('E50').Value = '  -4.24%  '

$ws.Range('D51').Value = '1.217'
$ws.Range('E51').Value = '  -1.92%  '
